$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates scraped from the upstream cryptos.xlsx refresh.
# Numeric-looking values in column D are prefixed with a leading apostrophe
# so Excel stores them as text (preserving trailing zeros / exact formatting)
# instead of silently converting them to numbers, matching the source data.
$ws.Range('D2').Value = '29.352.27'
$ws.Range('E2').Value = '  -0.50%  '
$ws.Range('D3').Value = '1.843.53'
$ws.Range('E3').Value = '  -0.79%  '
$ws.Range('D4').Value = '''0.9986'
$ws.Range('E4').Value = '  -0.48%  '
$ws.Range('D5').Value = '''240.23'
$ws.Range('E5').Value = '  -0.64%  '
$ws.Range('D6').Value = '''0.6298'
$ws.Range('E6').Value = '  -0.71%  '
$ws.Range('D7').Value = '''0.9999'
$ws.Range('D8').Value = '''0.07437'
$ws.Range('E8').Value = '  -1.89%  '
$ws.Range('D9').Value = '''0.2905'
$ws.Range('E9').Value = '  -0.82%  '
$ws.Range('E10').Value = '  +0.84%  '
$ws.Range('D11').Value = '''0.07739'
$ws.Range('E11').Value = '  -0.41%  '
$ws.Range('D12').Value = '1.847.66'
$ws.Range('E12').Value = '  -0.56%  '
$ws.Range('D13').Value = '''4.985'
$ws.Range('E13').Value = '  -1.20%  '
$ws.Range('D14').Value = '''0.6794'
$ws.Range('E14').Value = '  -1.03%  '
$ws.Range('D15').Value = '''0.00001023'
$ws.Range('E15').Value = '  -2.61%  '
$ws.Range('D16').Value = '''82.04'
$ws.Range('E16').Value = '  -1.76%  '
$ws.Range('D17').Value = '''6.268'
$ws.Range('E17').Value = '  +1.82%  '
$ws.Range('D18').Value = '29.346.99'
$ws.Range('E18').Value = '  -0.58%  '
$ws.Range('D19').Value = '''229.49'
$ws.Range('E19').Value = '  -0.39%  '
$ws.Range('E20').Value = '  -0.61%  '
$ws.Range('D22').Value = '''7.429'
$ws.Range('E22').Value = '  -1.14%  '
$ws.Range('E23').Value = '  -0.45%  '
$ws.Range('D24').Value = '''158.12'
$ws.Range('E24').Value = '  -0.93%  '
$ws.Range('E25').Value = '  +0.18%  '
$ws.Range('D26').Value = '''0.1354'
$ws.Range('E27').Value = '  -1.68%  '
$ws.Range('D28').Value = '''0.06528'
$ws.Range('E28').Value = '  +14.44%  '
$ws.Range('D29').Value = '''1.449'
$ws.Range('D30').Value = '''1.486'
$ws.Range('E30').Value = '  +0.31%  '
$ws.Range('E31').Value = '  -2.04%  '
$ws.Range('D32').Value = '''4.060'
$ws.Range('E32').Value = '  -0.31%  '
$ws.Range('E33').Value = '  +0.52%  '
$ws.Range('E34').Value = '  -1.78%  '
$ws.Range('D35').Value = '''0.6963'
$ws.Range('E35').Value = '  -0.20%  '
$ws.Range('D36').Value = '''2.572'
$ws.Range('E36').Value = '  -0.91%  '
$ws.Range('D37').Value = '''0.01853'
$ws.Range('E37').Value = '  +0.80%  '
$ws.Range('D38').Value = '''2.813'
$ws.Range('E38').Value = '  +0.94%  '
$ws.Range('D39').Value = '1.242.93'
$ws.Range('E39').Value = '  -1.17%  '
$ws.Range('D40').Value = '''6.799'
$ws.Range('E40').Value = '  +4.37%  '
$ws.Range('D41').Value = '''0.9339'
$ws.Range('E41').Value = '  +2.73%  '
$ws.Range('D42').Value = '''0.9995'
$ws.Range('E42').Value = '  -0.37%  '
$ws.Range('D43').Value = '1.991.39'
$ws.Range('E43').Value = '  -1.43%  '
$ws.Range('D44').Value = '''100.82'
$ws.Range('E44').Value = '  -0.83%  '
$ws.Range('E45').Value = '  -1.00%  '
$ws.Range('E46').Value = '  +3.39%  '
$ws.Range('D47').Value = '''7.052'
$ws.Range('E47').Value = '  -1.28%  '
$ws.Range('D48').Value = '''1.712'
$ws.Range('E48').Value = '  +1.93%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = '''9.017'
$ws.Range('E49').Value = '  -0.72%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').Value = '''0.1150'
$ws.Range('E50').Value = '  -1.81%  '
$ws.Range('D51').Value = '''0.3895'
$ws.Range('E51').Value = '  -2.23%  '
